# Add a new "2022-Q1" sheet (copied formatting from the "2021-Q4" sheet)
# positioned right before the "总计" (totals) sheet, and add a matching
# summary row at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet before "总计"
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Copy the cell formatting (styles/borders/fonts) from the 2021-Q4 sheet so
# the new sheet matches the look of the other quarterly sheets.
$template.Range("A1:H2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Index column
$newSheet.Range("A2").Value = 0

# Data row - write the text-like values through a scratch range that is
# formatted as text, so strings such as "003745" / "37.08" are preserved
# exactly (instead of being auto-converted into numbers).
$scratch = $newSheet.Range("Z1:AE1")
$scratch.NumberFormat = "@"
$newSheet.Range("Z1").Value = "003745"
$newSheet.Range("AA1").Value = "广发多元新兴股票"
$newSheet.Range("AB1").Value = "37.08"
$newSheet.Range("AC1").Value = "90.97"
$newSheet.Range("AD1").Value = "4.37"
$newSheet.Range("AE1").Value = "1.6204"

$scratch.Copy()
$newSheet.Range("B2").PasteSpecial(-4163)
$scratch.Clear()

$newSheet.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2. Insert a new top data row in the "总计" sheet for 2022-Q1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")

$ws.Rows.Item(2).Insert()
$ws.Range("A2:D2").ClearFormats()

# Re-use the formatting of the (now shifted) row below for the index cell.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1.62

# Renumber the index column for the rows that shifted down.
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
